$wb = $excel.ActiveWorkbook

# --- Sheet "Recommandations" (sheet1) ---
$ws1 = $wb.Worksheets.Item("Recommandations")

$ws1.Range("D2").Value = 3445.18
$ws1.Range("E2").Value = 112.96
$ws1.Range("D3").Value = 2850
$ws1.Range("E3").Value = 970
$ws1.Range("E5").Value = 680
$ws1.Range("D6").Value = 2650.53
$ws1.Range("E6").Value = 662.63
$ws1.Range("D7").Value = 2380
$ws1.Range("E7").Value = 590
$ws1.Range("D8").Value = 2365
$ws1.Range("E9").Value = 540
$ws1.Range("D10").Value = 2135
$ws1.Range("D11").Value = 1574.44
$ws1.Range("E11").Value = 382.27
$ws1.Range("D12").Value = 1490.89
$ws1.Range("E12").Value = 371.94
$ws1.Range("D13").Value = 1333.85
$ws1.Range("E13").Value = 336.77
$ws1.Range("D14").Value = 774.66
$ws1.Range("E14").Value = 194.07
$ws1.Range("D15").Value = 715.0599999999999
$ws1.Range("E15").Value = 178.09
$ws1.Range("D16").Value = 679.28
$ws1.Range("E16").Value = 171.22
$ws1.Range("D17").Value = 596.25
$ws1.Range("E17").Value = 143.45
$ws1.Range("D18").Value = 523.34
$ws1.Range("E18").Value = 131.54
$ws1.Range("D19").Value = 493.65
$ws1.Range("E19").Value = 123.32
$ws1.Range("D20").Value = 485.15
$ws1.Range("E20").Value = 121.2
$ws1.Range("D21").Value = 445.25
$ws1.Range("E21").Value = 111.06
$ws1.Range("D22").Value = 431.25
$ws1.Range("E22").Value = 107.85
$ws1.Range("D23").Value = 390.37
$ws1.Range("E23").Value = 97.44
$ws1.Range("A24").Value = 'UNILEVER CI (UNLC)'
$ws1.Range("D24").Value = 11.06
$ws1.Range("E24").Value = 7.47
$ws1.Range("A25").Value = 'ECOBANK COTE D''''IVOIRE (ECOC)'
$ws1.Range("B25").Value = 2
$ws1.Range("D25").Value = 7.44
$ws1.Range("E25").Value = 4.3
$ws1.Range("A26").Value = 'SOGB CI (SOGC)'
$ws1.Range("B26").Value = 1
$ws1.Range("D26").Value = 6.45
$ws1.Range("E26").Value = 6.45
$ws1.Range("A27").Value = 'ECOBANK TRANS. INCORP. TG (ETIT)'
$ws1.Range("D27").Value = 5.88
$ws1.Range("E27").Value = 5.88
$ws1.Range("A28").Value = 'CFAO MOTORS CI (CFAC)'
$ws1.Range("B28").Value = 2
$ws1.Range("C28").Value = 1
$ws1.Range("D28").Value = 5.44
$ws1.Range("E28").Value = 3.03
$ws1.Range("G28").Value = '👀 À surveiller'
$ws1.Range("A29").Value = 'BERNABE CI (BNBC)'
$ws1.Range("B29").Value = 2
$ws1.Range("C29").Value = 1
$ws1.Range("D29").Value = 5.35
$ws1.Range("E29").Value = 3.02
$ws1.Range("G29").Value = '👀 À surveiller'
$ws1.Range("A30").Value = 'NSIA BANQUE COTE D''IVOIRE (NSBC)'
$ws1.Range("B30").Value = 1
$ws1.Range("C30").Value = 0
$ws1.Range("D30").Value = 4.63
$ws1.Range("E30").Value = 4.63
$ws1.Range("G30").Value = '➖ Neutre'
$ws1.Range("A31").Value = 'SMB CI (SMBC)'
$ws1.Range("B31").Value = 1
$ws1.Range("C31").Value = 0
$ws1.Range("D31").Value = 3.74
$ws1.Range("E31").Value = 3.74
$ws1.Range("G31").Value = '➖ Neutre'
$ws1.Range("A32").Value = 'SUCRIVOIRE (SCRC)'
$ws1.Range("D32").Value = 3.09
$ws1.Range("E32").Value = 3.09
$ws1.Range("A33").Value = 'BANK OF AFRICA BN (BOAB)'
$ws1.Range("D33").Value = 2.81
$ws1.Range("E33").Value = 2.81
$ws1.Range("A34").Value = 'SETAO CI (STAC)'
$ws1.Range("C34").Value = 1
$ws1.Range("D34").Value = 2.18
$ws1.Range("E34").Value = -3.7
$ws1.Range("G34").Value = '👀 À surveiller'
$ws1.Range("A35").Value = 'SAFCA CI (SAFC)'
$ws1.Range("D35").Value = 0.2
$ws1.Range("E35").Value = 4.55
$ws1.Range("A36").Value = 'TOTAL'
$ws1.Range("B36").Value = 0
$ws1.Range("C36").Value = 4
$ws1.Range("D36").Value = 0
$ws1.Range("E36").Value = 0
$ws1.Range("G36").Value = '➖ Neutre'
$ws1.Range("A37").Value = 'SERVAIR ABIDJAN CI (ABJC)'
$ws1.Range("D37").Value = -0.01
$ws1.Range("E37").Value = -7.47
$ws1.Range("A38").Value = 'TOTALENERGIES MARKETING SN (TTLS)'
$ws1.Range("D38").Value = -0.62
$ws1.Range("E38").Value = 3.02
$ws1.Range("A39").Value = 'ORAGROUP TOGO (ORGT)'
$ws1.Range("C39").Value = 1
$ws1.Range("D39").Value = -2.42
$ws1.Range("E39").Value = -2.42
$ws1.Range("A40").Value = 'NEI-CEDA CI (NEIC)'
$ws1.Range("B40").Value = 0
$ws1.Range("D40").Value = -2.5
$ws1.Range("E40").Value = -2.5
$ws1.Range("G40").Value = '➖ Neutre'
$ws1.Range("A41").Value = 'BANK OF AFRICA CI (BOAC)'
$ws1.Range("D41").Value = -2.78
$ws1.Range("E41").Value = -2.78
$ws1.Range("A42").Value = 'VIVO ENERGY CI (SHEC)'
$ws1.Range("D42").Value = -3.21
$ws1.Range("E42").Value = -3.21
$ws1.Range("A43").Value = 'AFRICA GLOBAL LOGISTICS CI (SDSC)'
$ws1.Range("B43").Value = 1
$ws1.Range("C43").Value = 2
$ws1.Range("D43").Value = -4.16
$ws1.Range("E43").Value = -3.99
$ws1.Range("G43").Value = '👀 À surveiller'
$ws1.Range("A44").Value = 'TRACTAFRIC MOTORS CI (PRSC)'
$ws1.Range("D44").Value = -4.18
$ws1.Range("E44").Value = -4.18
$ws1.Range("A45").Value = 'ONATEL BF (ONTBF)'
$ws1.Range("B45").Value = 1
$ws1.Range("D45").Value = -5.05
$ws1.Range("E45").Value = 7.48
$ws1.Range("G45").Value = '👀 À surveiller'
$ws1.Range("A46").Value = 'BANK OF AFRICA BF (BOABF)'
$ws1.Range("C46").Value = 1
$ws1.Range("D46").Value = -5.33
$ws1.Range("E46").Value = -5.33
$ws1.Range("A48").Value = 'FILTISAC CI (FTSC)'
$ws1.Range("C48").Value = 3
$ws1.Range("D48").Value = -11.49
$ws1.Range("E48").Value = -5.66
$ws1.Range("F48").Value = '🔴 Vente'
$ws1.Range("G48").Value = '⚠️ Risque de décrochage'

# Row 49 (ONATEL BF) no longer present -> delete the row, shifting 2:49 -> dimension A1:G48
$ws1.Rows.Item(49).Delete()

# --- Sheet "Top_YTD" (sheet2) ---
$ws2 = $wb.Worksheets.Item("Top_YTD")

$ws2.Range("B2").Value = 10525576.66
$ws2.Range("B5").Value = 338154.99
$ws2.Range("B6").Value = 233189
$ws2.Range("B7").Value = 228177.96
$ws2.Range("B9").Value = 161057.12
$ws2.Range("B10").Value = 115628.52
$ws2.Range("B11").Value = 59234.27

Write-Host "Update complete"